$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update customer name and phone number in row 4
$ws.Range("B4").Value = "BS.Hoàng Văn Lợi"
$ws.Range("C4").Value = "'0978441444"

# Xóa dòng cuối cùng (rows 6 and 7)
$ws.Range("A6:G7").Clear()

$ws.Range("C5").Select()
